# Fruta / hortaliza, semanal
# Insert a new weekly record row above the existing row 8 (date 2021-10-22 /
# serial 44491), pushing that row down to row 9 unchanged, then populate the
# new row 8 with the latest weekly observation (date 2021-11-09 / serial
# 44509) using the same market/product metadata.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing row 8 (and everything after it, i.e. none) down to row 9,
# preserving its values/formatting intact.
$ws.Rows("8:8").Insert()

# Populate the newly inserted row 8 with the new weekly entry.
$ws.Cells.Item(8, 1).Value = 4
$ws.Cells.Item(8, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(8, 3).Value = "Los Lagos"
$ws.Cells.Item(8, 4).Value = 44509
$ws.Cells.Item(8, 5).Value = 10
$ws.Cells.Item(8, 6).Value = "Fruta"
$ws.Cells.Item(8, 7).Value = 100107
$ws.Cells.Item(8, 8).Value = "Otros"
$ws.Cells.Item(8, 9).Value = 100107002
$ws.Cells.Item(8, 10).Value = "Chirimoya"
$ws.Cells.Item(8, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(8, 12).Value = "Segunda"
$ws.Cells.Item(8, 13).Value = 200
$ws.Cells.Item(8, 14).Value = 19000
$ws.Cells.Item(8, 15).Value = 20000
$ws.Cells.Item(8, 16).Value = 19500
$ws.Cells.Item(8, 17).Value = "$/bandeja 8 kilos"
$ws.Cells.Item(8, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(8, 19).Value = 2438
$ws.Cells.Item(8, 20).Value = 8
